# Fruta / hortaliza, semanal
# Insert two new weekly price rows above the existing row 650 block,
# pushing all subsequent rows (650-722) down to (652-724), then
# populate the two newly inserted rows with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before row 650 (shifts old 650:722 -> 652:724)
$ws.Rows("650:651").Insert()

# --- New row 650 ---
$ws.Cells.Item(650, 1).Value2 = 4
$ws.Cells.Item(650, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(650, 3).Value2 = "Los Lagos"
$ws.Cells.Item(650, 4).Value2 = 44946
$ws.Cells.Item(650, 5).Value2 = 10
$ws.Cells.Item(650, 6).Value2 = "Fruta"
$ws.Cells.Item(650, 7).Value2 = 100108
$ws.Cells.Item(650, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(650, 9).Value2 = 100108006
$ws.Cells.Item(650, 10).Value2 = "Plátano"
$ws.Cells.Item(650, 11).Value2 = "Sin especificar"
$ws.Cells.Item(650, 12).Value2 = "Pintón"
$ws.Cells.Item(650, 13).Value2 = 600
$ws.Cells.Item(650, 14).Value2 = 26000
$ws.Cells.Item(650, 15).Value2 = 26000
$ws.Cells.Item(650, 16).Value2 = 26000
$ws.Cells.Item(650, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(650, 18).Value2 = "Ecuador"
$ws.Cells.Item(650, 19).Value2 = 1300
$ws.Cells.Item(650, 20).Value2 = 20

# --- New row 651 ---
$ws.Cells.Item(651, 1).Value2 = 4
$ws.Cells.Item(651, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(651, 3).Value2 = "Los Lagos"
$ws.Cells.Item(651, 4).Value2 = 44946
$ws.Cells.Item(651, 5).Value2 = 10
$ws.Cells.Item(651, 6).Value2 = "Fruta"
$ws.Cells.Item(651, 7).Value2 = 100108
$ws.Cells.Item(651, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(651, 9).Value2 = 100108006
$ws.Cells.Item(651, 10).Value2 = "Plátano"
$ws.Cells.Item(651, 11).Value2 = "Sin especificar"
$ws.Cells.Item(651, 12).Value2 = "Primera Pintón"
$ws.Cells.Item(651, 13).Value2 = 1200
$ws.Cells.Item(651, 14).Value2 = 28000
$ws.Cells.Item(651, 15).Value2 = 29000
$ws.Cells.Item(651, 16).Value2 = 28500
$ws.Cells.Item(651, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(651, 18).Value2 = "Ecuador"
$ws.Cells.Item(651, 19).Value2 = 1425
$ws.Cells.Item(651, 20).Value2 = 20
